# "Fechas de Tareas.xlsx" - nueva edicion de las propiedades a desarrollar
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Values: clear A2:A5 first so the shared-string table re-allocates
#    the same way the original author's edit did (verified empirically),
#    then retype everything in the exact original authoring order.
# ---------------------------------------------------------------------
$ws.Range("A2:A5").ClearContents()
$ws.Range("A2").Value = " 1. Crear input comentario"
$ws.Range("A3").Value = " 2. Guardar texto del comentario en la lista de comentarios"
$ws.Range("A4").Value = " 3. Agregar imagen"
$ws.Range("A5").Value = " 4. Recortar imagen seleccionada"
$ws.Range("A6").Value = " 4.1.  Agregar opción para recortar la imagen seleccionada"
$ws.Range("H1").Value = "TERMINADO"
$ws.Range("I1").Value = "EN PROCESO"
$ws.Range("A8").Value = " 4.2. Agregar la imagen segun la edición de recorte echa"
$ws.Range("A7").Value = " 4.3. Agregar modal que selecciona y edite la imagen"
$ws.Range("A9").Value = " 5. Agregar estilos"

# Dates (cols B & C rows 2-5)
$ws.Range("B2").Value = 45527
$ws.Range("C2").Value = 45527
$ws.Range("B3").Value = 45530
$ws.Range("C3").Value = 45530
$ws.Range("B4").Value = 45530
$ws.Range("C4").Value = 45530
$ws.Range("B5").Value = 45530
$ws.Range("C5").Value = 45531

Write-Output "values done"

# ---------------------------------------------------------------------
# 2) Number formats for the date columns (B:C rows 2-5) -> built-in "d-mmm"
# ---------------------------------------------------------------------
$ws.Range("B2:C5").NumberFormat = "d-mmm"

Write-Output "numberformat done"

# ---------------------------------------------------------------------
# 3) Fills (themed, tinted) for the header pills (H1:I1) and the
#    "done" / "in progress" rows (row2 = green, row3 = gold)
# ---------------------------------------------------------------------
$ws.Range("H1").Interior.ThemeColor = 9
$ws.Range("H1").Interior.TintAndShade = 0.79998168889431442

$ws.Range("I1").Interior.ThemeColor = 7
$ws.Range("I1").Interior.TintAndShade = 0.79998168889431442

$ws.Range("A2:C2").Interior.ThemeColor = 9
$ws.Range("A2:C2").Interior.TintAndShade = 0.79998168889431442

$ws.Range("A3:C3").Interior.ThemeColor = 7
$ws.Range("A3:C3").Interior.TintAndShade = 0.79998168889431442

Write-Output "fills done"

